$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# August 2023 data (rows 104-134)
# Row 104
$ws.Cells.Item(104, 1).Value = 45139
$ws.Cells.Item(104, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(104, 4).Value = 6500
$ws.Cells.Item(104, 6).Value = 5000
$ws.Cells.Item(104, 9).Value = 10500
$ws.Cells.Item(104, 10).Value = 1000
$ws.Cells.Item(104, 18).Value = 1000
$ws.Cells.Item(104, 19).Value = 18000
$ws.Cells.Item(104, 22).Value = 1000

# Row 105
$ws.Cells.Item(105, 1).Value = 45140
$ws.Cells.Item(105, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(105, 2).Value = 100000
$ws.Cells.Item(105, 3).Value = 5000
$ws.Cells.Item(105, 4).Value = 204500
$ws.Cells.Item(105, 6).Value = 272900
$ws.Cells.Item(105, 11).Value = 110000
$ws.Cells.Item(105, 18).Value = 9500
$ws.Cells.Item(105, 19).Value = 1600
$ws.Cells.Item(105, 21).Value = 120000

# Row 106
$ws.Cells.Item(106, 1).Value = 45141
$ws.Cells.Item(106, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(106, 3).Value = 92000
$ws.Cells.Item(106, 4).Value = 177050
$ws.Cells.Item(106, 5).Value = 39550
$ws.Cells.Item(106, 6).Value = 23500
$ws.Cells.Item(106, 9).Value = 5250
$ws.Cells.Item(106, 10).Value = 1000
$ws.Cells.Item(106, 11).Value = 50000
$ws.Cells.Item(106, 15).Value = 5000
$ws.Cells.Item(106, 18).Value = 4000
$ws.Cells.Item(106, 19).Value = 9000
$ws.Cells.Item(106, 20).Value = 1800

# Row 107
$ws.Cells.Item(107, 1).Value = 45142
$ws.Cells.Item(107, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(107, 4).Value = 83975
$ws.Cells.Item(107, 5).Value = 63025
$ws.Cells.Item(107, 6).Value = 226400
$ws.Cells.Item(107, 9).Value = 3500
$ws.Cells.Item(107, 10).Value = 1000
$ws.Cells.Item(107, 13).Value = 75000
$ws.Cells.Item(107, 18).Value = 24100
$ws.Cells.Item(107, 20).Value = 1000

# Row 108
$ws.Cells.Item(108, 1).Value = 45143
$ws.Cells.Item(108, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(108, 4).Value = 107350
$ws.Cells.Item(108, 6).Value = 297200
$ws.Cells.Item(108, 9).Value = 4500
$ws.Cells.Item(108, 10).Value = 1000
$ws.Cells.Item(108, 14).Value = 9000
$ws.Cells.Item(108, 18).Value = 3600
$ws.Cells.Item(108, 21).Value = 25000
$ws.Cells.Item(108, 22).Value = 5000

# Row 109
$ws.Cells.Item(109, 1).Value = 45144
$ws.Cells.Item(109, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(109, 3).Value = 1000
$ws.Cells.Item(109, 4).Value = 73750
$ws.Cells.Item(109, 6).Value = 33400
$ws.Cells.Item(109, 7).Value = 1890450
$ws.Cells.Item(109, 9).Value = 3500
$ws.Cells.Item(109, 10).Value = 1000
$ws.Cells.Item(109, 18).Value = 9000
$ws.Cells.Item(109, 19).Value = 200

# Row 110
$ws.Cells.Item(110, 1).Value = 45145
$ws.Cells.Item(110, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(110, 3).Value = 92000
$ws.Cells.Item(110, 4).Value = 191950
$ws.Cells.Item(110, 9).Value = 2000
$ws.Cells.Item(110, 10).Value = 1000
$ws.Cells.Item(110, 17).Value = 3000
$ws.Cells.Item(110, 18).Value = 500
$ws.Cells.Item(110, 21).Value = 20000
$ws.Cells.Item(110, 22).Value = 125000

# Row 111
$ws.Cells.Item(111, 1).Value = 45146
$ws.Cells.Item(111, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(111, 3).Value = 1500
$ws.Cells.Item(111, 4).Value = 64000
$ws.Cells.Item(111, 6).Value = 34400
$ws.Cells.Item(111, 9).Value = 7000
$ws.Cells.Item(111, 10).Value = 1000
$ws.Cells.Item(111, 18).Value = 5000
$ws.Cells.Item(111, 22).Value = 1000
$ws.Cells.Item(111, 23).Value = 10000

# Row 112
$ws.Cells.Item(112, 1).Value = 45147
$ws.Cells.Item(112, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(112, 4).Value = 168500
$ws.Cells.Item(112, 6).Value = 20500
$ws.Cells.Item(112, 9).Value = 7000
$ws.Cells.Item(112, 10).Value = 1000
$ws.Cells.Item(112, 11).Value = 80000
$ws.Cells.Item(112, 18).Value = 28000
$ws.Cells.Item(112, 21).Value = 16000
$ws.Cells.Item(112, 22).Value = 397000
$ws.Cells.Item(112, 23).Value = 46500

# Row 113
$ws.Cells.Item(113, 1).Value = 45148
$ws.Cells.Item(113, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(113, 3).Value = 92000
$ws.Cells.Item(113, 4).Value = 71450
$ws.Cells.Item(113, 6).Value = 272000
$ws.Cells.Item(113, 10).Value = 1000
$ws.Cells.Item(113, 11).Value = 40000
$ws.Cells.Item(113, 13).Value = 300000
$ws.Cells.Item(113, 15).Value = 25000
$ws.Cells.Item(113, 18).Value = 42000
$ws.Cells.Item(113, 21).Value = 130000
$ws.Cells.Item(113, 22).Value = 80000
$ws.Cells.Item(113, 23).Value = 15000

# Row 114
$ws.Cells.Item(114, 1).Value = 45149
$ws.Cells.Item(114, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(114, 4).Value = 69600
$ws.Cells.Item(114, 6).Value = 417500
$ws.Cells.Item(114, 9).Value = 10500
$ws.Cells.Item(114, 10).Value = 1000
$ws.Cells.Item(114, 18).Value = 6500
$ws.Cells.Item(114, 20).Value = 1500
$ws.Cells.Item(114, 21).Value = 50000

# Row 115
$ws.Cells.Item(115, 1).Value = 45150
$ws.Cells.Item(115, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(115, 4).Value = 122500
$ws.Cells.Item(115, 6).Value = 92600
$ws.Cells.Item(115, 9).Value = 3500
$ws.Cells.Item(115, 10).Value = 1000
$ws.Cells.Item(115, 18).Value = 14500
$ws.Cells.Item(115, 21).Value = 39000
$ws.Cells.Item(115, 22).Value = 135000
$ws.Cells.Item(115, 23).Value = 25000

# Row 116
$ws.Cells.Item(116, 1).Value = 45151
$ws.Cells.Item(116, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(116, 3).Value = 92000
$ws.Cells.Item(116, 4).Value = 159375
$ws.Cells.Item(116, 6).Value = 529400
$ws.Cells.Item(116, 7).Value = 1726500
$ws.Cells.Item(116, 9).Value = 9000
$ws.Cells.Item(116, 10).Value = 1000

# Row 117
$ws.Cells.Item(117, 1).Value = 45152
$ws.Cells.Item(117, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(117, 4).Value = 140175
$ws.Cells.Item(117, 6).Value = 49100
$ws.Cells.Item(117, 10).Value = 1000
$ws.Cells.Item(117, 14).Value = 1000
$ws.Cells.Item(117, 18).Value = 4500
$ws.Cells.Item(117, 22).Value = 51000

# Row 118
$ws.Cells.Item(118, 1).Value = 45153
$ws.Cells.Item(118, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(118, 4).Value = 16500
$ws.Cells.Item(118, 6).Value = 4800
$ws.Cells.Item(118, 9).Value = 10500
$ws.Cells.Item(118, 10).Value = 1000
$ws.Cells.Item(118, 18).Value = 10000

# Row 119
$ws.Cells.Item(119, 1).Value = 45154
$ws.Cells.Item(119, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(119, 4).Value = 110250
$ws.Cells.Item(119, 6).Value = 961400
$ws.Cells.Item(119, 9).Value = 10500
$ws.Cells.Item(119, 10).Value = 1000
$ws.Cells.Item(119, 15).Value = 55000
$ws.Cells.Item(119, 18).Value = 56000
$ws.Cells.Item(119, 19).Value = 10000
$ws.Cells.Item(119, 21).Value = 40000
$ws.Cells.Item(119, 22).Value = 90000

# Row 120
$ws.Cells.Item(120, 1).Value = 45155
$ws.Cells.Item(120, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(120, 3).Value = 92000
$ws.Cells.Item(120, 4).Value = 155600
$ws.Cells.Item(120, 6).Value = 33900
$ws.Cells.Item(120, 9).Value = 10500
$ws.Cells.Item(120, 10).Value = 1000
$ws.Cells.Item(120, 11).Value = 30000
$ws.Cells.Item(120, 17).Value = 2000
$ws.Cells.Item(120, 18).Value = 9000
$ws.Cells.Item(120, 20).Value = 1500
$ws.Cells.Item(120, 22).Value = 338800
$ws.Cells.Item(120, 23).Value = 30000

# Row 121
$ws.Cells.Item(121, 1).Value = 45156
$ws.Cells.Item(121, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(121, 4).Value = 128700
$ws.Cells.Item(121, 9).Value = 10500
$ws.Cells.Item(121, 10).Value = 1000
$ws.Cells.Item(121, 18).Value = 28500
$ws.Cells.Item(121, 19).Value = 1000
$ws.Cells.Item(121, 21).Value = 15000

# Row 122
$ws.Cells.Item(122, 1).Value = 45157
$ws.Cells.Item(122, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(122, 4).Value = 100550
$ws.Cells.Item(122, 9).Value = 7000
$ws.Cells.Item(122, 10).Value = 1000
$ws.Cells.Item(122, 18).Value = 3000
$ws.Cells.Item(122, 20).Value = 1500
$ws.Cells.Item(122, 23).Value = 5000

# Row 123
$ws.Cells.Item(123, 1).Value = 45158
$ws.Cells.Item(123, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(123, 3).Value = 77500
$ws.Cells.Item(123, 4).Value = 102225
$ws.Cells.Item(123, 7).Value = 2060800
$ws.Cells.Item(123, 9).Value = 10500
$ws.Cells.Item(123, 20).Value = 1500

# Row 124
$ws.Cells.Item(124, 1).Value = 45159
$ws.Cells.Item(124, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(124, 4).Value = 97750
$ws.Cells.Item(124, 5).Value = 37250
$ws.Cells.Item(124, 10).Value = 1000
$ws.Cells.Item(124, 18).Value = 9000
$ws.Cells.Item(124, 22).Value = 765000
$ws.Cells.Item(124, 23).Value = 185000

# Row 125
$ws.Cells.Item(125, 1).Value = 45160
$ws.Cells.Item(125, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(125, 4).Value = 13000
$ws.Cells.Item(125, 6).Value = 240900
$ws.Cells.Item(125, 9).Value = 10500
$ws.Cells.Item(125, 10).Value = 1000

# Row 126
$ws.Cells.Item(126, 1).Value = 45161
$ws.Cells.Item(126, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(126, 3).Value = 77500
$ws.Cells.Item(126, 4).Value = 132150
$ws.Cells.Item(126, 6).Value = 518700
$ws.Cells.Item(126, 9).Value = 10500
$ws.Cells.Item(126, 10).Value = 1000
$ws.Cells.Item(126, 11).Value = 160000
$ws.Cells.Item(126, 17).Value = 3000
$ws.Cells.Item(126, 18).Value = 30500
$ws.Cells.Item(126, 20).Value = 2100
$ws.Cells.Item(126, 21).Value = 145650

# Row 127
$ws.Cells.Item(127, 1).Value = 45162
$ws.Cells.Item(127, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(127, 4).Value = 119400
$ws.Cells.Item(127, 6).Value = 105600
$ws.Cells.Item(127, 10).Value = 1000
$ws.Cells.Item(127, 11).Value = 30000
$ws.Cells.Item(127, 18).Value = 9000
$ws.Cells.Item(127, 22).Value = 51750
$ws.Cells.Item(127, 23).Value = 57500

# Row 128
$ws.Cells.Item(128, 1).Value = 45163
$ws.Cells.Item(128, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(128, 4).Value = 96900
$ws.Cells.Item(128, 6).Value = 86000
$ws.Cells.Item(128, 9).Value = 21000
$ws.Cells.Item(128, 23).Value = 15000

# Row 129
$ws.Cells.Item(129, 1).Value = 45164
$ws.Cells.Item(129, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(129, 3).Value = 77500
$ws.Cells.Item(129, 4).Value = 118225
$ws.Cells.Item(129, 6).Value = 545300
$ws.Cells.Item(129, 9).Value = 6500
$ws.Cells.Item(129, 10).Value = 1000
$ws.Cells.Item(129, 18).Value = 5500
$ws.Cells.Item(129, 23).Value = 26000

# Row 130
$ws.Cells.Item(130, 1).Value = 45165
$ws.Cells.Item(130, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(130, 4).Value = 95650
$ws.Cells.Item(130, 6).Value = 103000
$ws.Cells.Item(130, 7).Value = 1643150
$ws.Cells.Item(130, 9).Value = 10500
$ws.Cells.Item(130, 10).Value = 1000
$ws.Cells.Item(130, 13).Value = 30000
$ws.Cells.Item(130, 18).Value = 9000

# Row 131
$ws.Cells.Item(131, 1).Value = 45166
$ws.Cells.Item(131, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(131, 3).Value = 77500
$ws.Cells.Item(131, 4).Value = 223500
$ws.Cells.Item(131, 6).Value = 2500
$ws.Cells.Item(131, 9).Value = 2000
$ws.Cells.Item(131, 10).Value = 500
$ws.Cells.Item(131, 13).Value = 300000
$ws.Cells.Item(131, 18).Value = 500
$ws.Cells.Item(131, 19).Value = 500

# Row 132
$ws.Cells.Item(132, 1).Value = 45167
$ws.Cells.Item(132, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(132, 4).Value = 35000
$ws.Cells.Item(132, 6).Value = 29000
$ws.Cells.Item(132, 9).Value = 10500
$ws.Cells.Item(132, 15).Value = 5000
$ws.Cells.Item(132, 18).Value = 9000
$ws.Cells.Item(132, 22).Value = 17500

# Row 133
$ws.Cells.Item(133, 1).Value = 45168
$ws.Cells.Item(133, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(133, 4).Value = 101550
$ws.Cells.Item(133, 6).Value = 1395200
$ws.Cells.Item(133, 9).Value = 10500
$ws.Cells.Item(133, 10).Value = 1000
$ws.Cells.Item(133, 11).Value = 60000
$ws.Cells.Item(133, 18).Value = 105500
$ws.Cells.Item(133, 19).Value = 3000
$ws.Cells.Item(133, 22).Value = 15000

# Row 134
$ws.Cells.Item(134, 1).Value = 45169
$ws.Cells.Item(134, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(134, 3).Value = 77500
$ws.Cells.Item(134, 4).Value = 69950
$ws.Cells.Item(134, 6).Value = 42500
$ws.Cells.Item(134, 10).Value = 1000
$ws.Cells.Item(134, 11).Value = 30000
$ws.Cells.Item(134, 16).Value = 200000
$ws.Cells.Item(134, 18).Value = 6500

